# Add a new "Temperatures" worksheet after the existing sheets, populate it
# with a small time/temperature table, and make it the active sheet
# (matching the posted-after-class materials described in the commit).

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing worksheet so it lands at the
# end of the tab strip (Sheet1, Sheet2, Sheet3, Temperatures).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Temperatures"

# Header row.
$ws.Range("E3").Value = "Time"
$ws.Range("F3").Value = "Temp"

# Data rows.
$ws.Range("E4").Value = 12
$ws.Range("F4").Value = 98

$ws.Range("E5").Value = 13
$ws.Range("F5").Value = 99

$ws.Range("E6").Value = 14
$ws.Range("F6").Value = 97

# Select the populated range and make this the active/visible tab.
[void]$ws.Range("E3:F6").Select()
